$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.450.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "'1.640.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.35%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'1.003"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "'304.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").Value = "'0.3775"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("D8").Value = "'52.26"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.01%  "
$ws.Range("D9").Value = "'0.3642"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("D10").Value = "'1.247"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("D11").Value = "'0.08109"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "'22.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("D14").Value = "'6.641"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.95%  "
$ws.Range("D15").Value = "'0.00001253"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("D16").Value = "'7.291"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.57%  "
$ws.Range("D17").Value = "'1.634.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.79%  "
$ws.Range("D18").Value = "'94.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("D19").Value = "'0.06938"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("D20").Value = "'18.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").Value = "'6.547"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "'23.467.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.12%  "
$ws.Range("D24").Value = "'12.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").Value = "'3.222"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.96%  "
$ws.Range("D26").Value = "'2.453"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("D27").Value = "'21.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("D28").Value = "'150.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = "'5.311"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.99%  "
$ws.Range("D30").Value = "'135.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("D31").Value = "'2.317"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.81%  "
$ws.Range("D32").Value = "'1.816.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.09%  "
$ws.Range("D33").Value = "'6.884"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.65%  "
$ws.Range("D34").Value = "'10.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.07%  "
$ws.Range("D35").Value = "'0.9638"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.54%  "
$ws.Range("D36").Value = "'0.02859"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.52%  "
$ws.Range("D37").Value = "'6.269"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.08%  "
$ws.Range("D38").Value = "'0.2555"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.83%  "
$ws.Range("D39").Value = "'0.07284"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.64%  "
$ws.Range("D40").Value = "'0.08866"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.39%  "
$ws.Range("D41").Value = "'1.375"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.87%  "
$ws.Range("D42").Value = "'0.7125"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.61%  "
$ws.Range("D43").Value = "'16.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.96%  "
$ws.Range("D44").Value = "'12.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.58%  "
$ws.Range("D45").Value = "'0.6561"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("D46").Value = "'2.356"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.54%  "
$ws.Range("D47").Value = "'1.002"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("D48").Value = "'3.997"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("D49").Value = "'0.07995"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'1.218"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.12%  "
$ws.Range("D51").Value = "'127.99"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.36%  "
